# Reversing hotfix on Issue608a: restore the "AccountData" box + its
# connector, and shift a whole row of boxes/connectors left to make room.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points while
# the OOXML we are targeting is specified in EMU (1 pt = 12700 EMU). The
# host's float plumbing truncates when it converts back to EMU, so a
# naive "emu / 12700.0" can land one EMU short. Nudging by half an EMU
# before the division keeps the truncation from eating a unit and gives
# us an exact EMU round-trip.
function EMUToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Move-Shape($shapes, $id, $x, $y, $cx, $cy) {
    $shape = Get-ShapeById $shapes $id
    if ($null -eq $shape) {
        return $null
    }
    if ($null -ne $x) { $shape.Left = EMUToPt $x }
    if ($null -ne $y) { $shape.Top = EMUToPt $y }
    if ($null -ne $cx) { $shape.Width = EMUToPt $cx }
    if ($null -ne $cy) { $shape.Height = EMUToPt $cy }
    return $shape
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# "BaseData" rectangle
Move-Shape $shapes 38 4114800 2560637 $null $null | Out-Null

# Elbow Connector 38 (horizontal bent connector under the data boxes)
Move-Shape $shapes 39 2819401 3475036 $null $null | Out-Null

# Elbow Connector 68 that points into the Isosceles Triangle (id 37, idx 3).
# It additionally gains a begin-connection to shape 174 (idx 0) and a
# horizontal flip in the source XML; the connection-site wiring isn't
# reachable through the exposed ConnectorFormat API in this host, so we
# still apply the geometry + flip changes that are controllable here.
$conn69 = Move-Shape $shapes 69 4891088 3194050 18916 966787
if ($null -ne $conn69) {
    $conn69.HorizontalFlip = -1
}

# Elbow Connector 68 (several plain straight connectors below the boxes)
Move-Shape $shapes 99 5715000 3475037 $null $null | Out-Null
Move-Shape $shapes 102 6477000 3475037 $null $null | Out-Null
Move-Shape $shapes 105 2819400 3475037 $null $null | Out-Null

# Isosceles Triangle 36 (arrow head)
Move-Shape $shapes 37 4752975 2965450 $null $null | Out-Null

# "StudentData" rectangle
Move-Shape $shapes 172 5029200 3703637 $null $null | Out-Null

# Elbow Connector 68
Move-Shape $shapes 108 3886200 3475037 $null $null | Out-Null

# "CourseData" rectangle
Move-Shape $shapes 179 3048000 3703637 $null $null | Out-Null

# "InstructorData" rectangle
Move-Shape $shapes 173 5943600 4160837 $null $null | Out-Null

# "EvaluationData" rectangle
Move-Shape $shapes 174 4114800 4160837 $null $null | Out-Null

# "SubmissionData" rectangle
Move-Shape $shapes 175 2286000 4160837 $null $null | Out-Null

# Restore the "AccountData" rectangle by cloning its sibling "StudentData"
# box (same fill/line/effect style refs, same text formatting) and then
# re-pointing it at the new position/size/text.
$studentData = Get-ShapeById $shapes 172
$accountData = $studentData.Duplicate().Item(1)
$accountData.Name = "Rectangle 47"
$accountData.TextFrame.TextRange.Text = "AccountData"
$accountData.Left = EMUToPt 6553200
$accountData.Top = EMUToPt 3703637
$accountData.Width = EMUToPt 1304746
$accountData.Height = EMUToPt 361770

# Restore the connector feeding the "AccountData" box by cloning one of
# the other plain elbow connectors (same straightConnector1 geometry and
# accent6 style) and repositioning it.
$siblingConn = Get-ShapeById $shapes 99
$accountConn = $siblingConn.Duplicate().Item(1)
$accountConn.Left = EMUToPt 7010400
$accountConn.Top = EMUToPt 3475037
$accountConn.Width = EMUToPt 0
$accountConn.Height = EMUToPt 228600
